# Atualização automática de TRES_COROAS.xlsx
# Rename worksheets:
#   "Paineis DARQ"               -> "PAINEIS DARQ"
#   "Recolhimento x Eliminacao"  -> "RECOLHIMENTO X ELIMINAÇÃO"

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"
